$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '64.669.90'
$ws.Range("E2").Value = '  -1.74%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.419.94'
$ws.Range("E3").Value = '  -2.45%  '
$ws.Range("E4").Value = '  -0.05%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '572.40'
$ws.Range("E5").Value = '  -1.50%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '158.58'
$ws.Range("E6").Value = '  -2.08%  '
$ws.Range("E7").Value = '  +0.02%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.589'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '3.422.91'
$ws.Range("E9").Value = '  -2.50%  '
$ws.Range("E10").Value = '  -1.94%  '
$ws.Range("E11").Value = '  -3.05%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.437'
$ws.Range("E12").Value = '  -2.19%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.011.63'
$ws.Range("E13").Value = '  -2.45%  '
$ws.Range("E14").Value = '  -0.45%  '
$ws.Range("E15").Value = '  -4.43%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '27.66'
$ws.Range("E16").Value = '  -3.75%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '64.730.00'
$ws.Range("E17").Value = '  -1.65%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.419.50'
$ws.Range("E18").Value = '  -2.70%  '
$ws.Range("E19").Value = '  -2.06%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '13.84'
$ws.Range("E20").Value = '  -3.47%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '379.72'
$ws.Range("E21").Value = '  -2.95%  '
$ws.Range("E22").Value = '  -3.78%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.547'
$ws.Range("E23").Value = '  -1.12%  '
$ws.Range("E24").Value = '  +0.06%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '71.84'
$ws.Range("E26").Value = '  -5.59%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.90'
$ws.Range("E27").Value = '  -0.75%  '
$ws.Range("E28").Value = '  -1.45%  '
$ws.Range("E29").Value = '  -0.29%  '
$ws.Range("E30").Value = '  +0.15%  '
$ws.Range("E31").Value = '  -3.93%  '
$ws.Range("E32").Value = '  -3.00%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '23.18'
$ws.Range("E33").Value = '  -2.56%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '7.01'
$ws.Range("E34").Value = '  -2.54%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.57'
$ws.Range("E35").Value = '  +0.19%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '160.81'
$ws.Range("E36").Value = '  -1.41%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.89'
$ws.Range("E37").Value = '  -3.35%  '
$ws.Range("B38").Value = 'Maker'
$ws.Range("C38").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.887.02'
$ws.Range("E38").Value = '  -6.58%  '
$ws.Range("B39").Value = 'Hedera'
$ws.Range("C39").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0750'
$ws.Range("E39").Value = '  -3.29%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '6.72'
$ws.Range("E40").Value = '  +3.26%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '26.23'
$ws.Range("E41").Value = '  -4.35%  '
$ws.Range("B42").Value = 'Filecoin'
$ws.Range("C42").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '4.53'
$ws.Range("E42").Value = '  -0.47%  '
$ws.Range("B43").Value = 'OKB'
$ws.Range("C43").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '43.09'
$ws.Range("E43").Value = '  -0.31%  '
$ws.Range("E44").Value = '  -2.71%  '
$ws.Range("E45").Value = '  -1.86%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '25.84'
$ws.Range("E46").Value = '  -0.35%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '317.82'
$ws.Range("E47").Value = '  +0.87%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.25'
$ws.Range("E48").Value = '  -0.81%  '
$ws.Range("E49").Value = '  -4.79%  '
$ws.Range("E50").Value = '  -3.29%  '
$ws.Range("E51").Value = '  -2.36%  '
